$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (pushes RJ-45 jack / 3D printing filament / 22 AWG
# hookup wire rows down by one) to make room for the missing "Diode" BOM entry.
$ws.Rows("14:14").Insert()

# Populate the new row with the diode part info.
$ws.Range("A14").Value = "Diode"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "1N4007"
$ws.Range("D14").Value = "https://octopart.com/1n4007-t-diodes+inc.-55389582?r=sp&s=AmIHlQKgSlimIGmMg_WMyg"

# The MOSFET row's purchase-URL cell (D13) previously had plain text with no live
# link; turn it into a real hyperlink (re-applying the built-in "Hyperlink" cell
# style afterwards so it keeps reusing the existing style record rather than a
# freshly-minted one).
$ws.Hyperlinks.Add($ws.Range("D13"), "https://octopart.com/search?r=sp&s=afxfcSvbQpS9Ntd6rMgB_Q&q=2N7000&start=0") | Out-Null
$ws.Range("D13").Style = "Hyperlink"

# The new diode row's purchase-URL cell (D14) just needs the same visual
# "Hyperlink" cell style applied (matching D13's look) without becoming a live
# hyperlink itself.
$ws.Range("D14").Style = "Hyperlink"

# Grow Table1 so the new row participates in the table/autofilter range.
$lo = $wb.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D17"))

# Match the selection left behind by the edit session.
$ws.Range("D14").Select() | Out-Null
